$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings are created in this specific order (first use) to
# --- mirror the author's actual edit order from the source diff.

# 1) "N/A" AI type for both Kitsune forms (row 17, 18)
$ws.Range("L17").Value = "N/A"
$ws.Range("L18").Value = "N/A"

# 2) "sword " (trailing space) - 3rd spell of Kitsunario's default group
$ws.Range("Y17").Value = "sword "

# 3) New "GROUP/HEALTH_CRITICAL" spell group - first introduced on Lilim (row 5)
$ws.Range("AC5").Value = "GROUP/HEALTH_CRITICAL"
$ws.Range("AI17").Value = "GROUP/HEALTH_CRITICAL"
$ws.Range("AI18").Value = "GROUP/HEALTH_CRITICAL"

# 4) "halberd" new weapon root
$ws.Range("AC17").Value = "halberd"
$ws.Range("AJ18").Value = "halberd"

# 5) "empowered" / "weakened" new spell styles
$ws.Range("AL17").Value = "empowered"
$ws.Range("AL18").Value = "weakened"

# 6) "FormChange1" new AI type used by the Tanuki forms that change into Kitsune
$ws.Range("L11").Value = "FormChange1"
$ws.Range("L12").Value = "FormChange1"

# 7) AI parameters describing the form change targets
$ws.Range("M11").Value = "HealthLow1,Kitsunario,false"
$ws.Range("M12").Value = "HealthLow1,Kitsunuigi,false"

# --- Row 5 (Lilim): extend with a new GROUP/HEALTH_CRITICAL spell ---
$ws.Range("Z5").Value = "lance "
$ws.Range("AA5").Value = "veld"
$ws.Range("AB5").Value = "null"
$ws.Range("AD5").Value = "sword"
$ws.Range("AE5").Value = "agni"
$ws.Range("AF5").Value = "null"
$ws.Range("AG5").Value = "END"

# --- Row 11 (Tanooki): stat/AI rebalance ---
$ws.Range("D11").Value = 100
$ws.Range("F11").Value = 3
$ws.Range("J11").Value = 0.9
$ws.Range("K11").Value = 10
$ws.Range("S11").Value = "sword"

# --- Row 12 (Tater Totnuki): stat/AI rebalance ---
$ws.Range("D12").Value = 150
$ws.Range("F12").Value = 2
$ws.Range("J12").Value = 0.9
$ws.Range("K12").Value = 10

# --- Row 17 (Kitsunario): stats + AI + full spell-group rebuild ---
$ws.Range("F17").Value = 3
$ws.Range("I17").Value = 0.75
$ws.Range("J17").Value = 0.9
$ws.Range("K17").Value = 10
$ws.Range("N17").Value = -2
$ws.Range("O17").Value = -2
$ws.Range("P17").Value = -2
$ws.Range("T17").Value = "veld"
$ws.Range("V17").Value = "lance"
$ws.Range("W17").Value = "cryo"
$ws.Range("X17").Value = "null"
$ws.Range("Z17").Value = "null"
$ws.Range("AA17").Value = "null"
$ws.Range("AB17").Value = "GROUP/HEALTH_LOW"
$ws.Range("AD17").Value = "null"
$ws.Range("AE17").Value = "aimed"
$ws.Range("AF17").Value = "sword"
$ws.Range("AG17").Value = "null"
$ws.Range("AH17").Value = "null"
$ws.Range("AJ17").Value = "selfcare"
$ws.Range("AK17").Value = "agni"
$ws.Range("AM17").Value = "END"

# --- Row 18 (Kitsunuigi): stats + AI + full spell-group rebuild ---
$ws.Range("D18").Value = 150
$ws.Range("F18").Value = 2
$ws.Range("I18").Value = 0.75
$ws.Range("J18").Value = 0.9
$ws.Range("K18").Value = 10
$ws.Range("N18").Value = -2
$ws.Range("O18").Value = 2
$ws.Range("P18").Value = -2
$ws.Range("Q18").Value = -2
$ws.Range("T18").Value = "agni"
$ws.Range("V18").Value = "sword"
$ws.Range("W18").Value = "null"
$ws.Range("X18").Value = "null"
$ws.Range("Y18").Value = "lance"
$ws.Range("Z18").Value = "veld"
$ws.Range("AA18").Value = "null"
$ws.Range("AB18").Value = "GROUP/HEALTH_LOW"
$ws.Range("AC18").Value = "lance"
$ws.Range("AD18").Value = "agni"
$ws.Range("AE18").Value = "null"
$ws.Range("AF18").Value = "sword"
$ws.Range("AG18").Value = "null"
$ws.Range("AH18").Value = "null"
$ws.Range("AK18").Value = "cryo"
$ws.Range("AM18").Value = "END"

# --- Cosmetic: column width + cursor selection ---
$ws.Columns.Item(12).ColumnWidth = 11.33
$ws.Range("M13").Select()
